$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 12.265976
$ws.Range("H2").Value = 36.797928
$ws.Range("I2").Value = 0.004000867643088759
$ws.Range("J2").Value = 0.004000867643088758
$ws.Range("M2").Value = 7.407905
$ws.Range("N2").Value = 22.223715
$ws.Range("O2").Value = 0.1577242380174723
$ws.Range("P2").Value = 0.1577242380174723
$ws.Range("Q2").Value = 90.86518494028
$ws.Range("R2").Value = 817.7866644625199
$ws.Range("S2").Value = 0.0006310338004149349
$ws.Range("T2").Value = 0.0006310338004149348
# Row 3
$ws.Range("G3").Value = 12.265976
$ws.Range("H3").Value = 36.797928
$ws.Range("I3").Value = 0.004000867643088759
$ws.Range("J3").Value = 0.004000867643088758
$ws.Range("O3").Value = 0.3510414535684271
$ws.Range("P3").Value = 0.3510414535684271
$ws.Range("Q3").Value = 202.2355409741547
$ws.Range("R3").Value = 1820.119868767392
$ws.Range("S3").Value = 0.001404470392964765
$ws.Range("T3").Value = 0.001404470392964765
# Row 4
$ws.Range("G4").Value = 12.265976
$ws.Range("H4").Value = 36.797928
$ws.Range("I4").Value = 0.004000867643088759
$ws.Range("J4").Value = 0.004000867643088758
$ws.Range("M4").Value = 5.464566666666666
$ws.Range("N4").Value = 16.3937
$ws.Range("O4").Value = 0.1163479571613943
$ws.Range("P4").Value = 0.1163479571613943
$ws.Range("Q4").Value = 67.02824358373333
$ws.Range("R4").Value = 603.2541922536
$ws.Range("S4").Value = 0.0004654927771464995
$ws.Range("T4").Value = 0.0004654927771464994
# Row 5
$ws.Range("G5").Value = 12.265976
$ws.Range("H5").Value = 36.797928
$ws.Range("I5").Value = 0.004000867643088759
$ws.Range("J5").Value = 0.004000867643088758
$ws.Range("M5").Value = 17.60745533333333
$ws.Range("N5").Value = 52.822366
$ws.Range("O5").Value = 0.3748863512527063
$ws.Range("P5").Value = 0.3748863512527063
$ws.Range("Q5").Value = 215.9726245397387
$ws.Range("R5").Value = 1943.753620857648
$ws.Range("S5").Value = 0.00149987067256256
$ws.Range("T5").Value = 0.001499870672562559
# Row 6
$ws.Range("I6").Value = 0.9924545876219728
$ws.Range("J6").Value = 0.9924545876219727
$ws.Range("M6").Value = 7.407905
$ws.Range("N6").Value = 22.223715
$ws.Range("O6").Value = 0.1577242380174723
$ws.Range("P6").Value = 0.1577242380174723
$ws.Range("Q6").Value = 22540.00324276642
$ws.Range("R6").Value = 202860.0291848978
$ws.Range("S6").Value = 0.1565341435996203
$ws.Range("T6").Value = 0.1565341435996204
# Row 7
$ws.Range("I7").Value = 0.9924545876219728
$ws.Range("J7").Value = 0.9924545876219727
$ws.Range("O7").Value = 0.3510414535684271
$ws.Range("P7").Value = 0.3510414535684271
$ws.Range("S7").Value = 0.3483927010394712
$ws.Range("T7").Value = 0.3483927010394712
# Row 8
$ws.Range("I8").Value = 0.9924545876219728
$ws.Range("J8").Value = 0.9924545876219727
$ws.Range("M8").Value = 5.464566666666666
$ws.Range("N8").Value = 16.3937
$ws.Range("O8").Value = 0.1163479571613943
$ws.Range("P8").Value = 0.1163479571613943
$ws.Range("Q8").Value = 16627.01538248398
$ws.Range("R8").Value = 149643.1384423558
$ws.Range("S8").Value = 0.1154700638452705
$ws.Range("T8").Value = 0.1154700638452705
# Row 9
$ws.Range("I9").Value = 0.9924545876219728
$ws.Range("J9").Value = 0.9924545876219727
$ws.Range("M9").Value = 17.60745533333333
$ws.Range("N9").Value = 52.822366
$ws.Range("O9").Value = 0.3748863512527063
$ws.Range("P9").Value = 0.3748863512527063
$ws.Range("Q9").Value = 53574.13469937834
$ws.Range("R9").Value = 482167.212294405
$ws.Range("S9").Value = 0.3720576791376106
$ws.Range("T9").Value = 0.3720576791376106
# Row 10
$ws.Range("G10").Value = 8.377189333333332
$ws.Range("H10").Value = 25.131568
$ws.Range("I10").Value = 0.002732438555542716
$ws.Range("J10").Value = 0.002732438555542716
$ws.Range("M10").Value = 7.407905
$ws.Range("N10").Value = 22.223715
$ws.Range("O10").Value = 0.1577242380174723
$ws.Range("P10").Value = 0.1577242380174723
$ws.Range("Q10").Value = 62.05742274834665
$ws.Range("R10").Value = 558.5168047351199
$ws.Range("S10").Value = 0.0004309717891025375
$ws.Range("T10").Value = 0.0004309717891025375
# Row 11
$ws.Range("G11").Value = 8.377189333333332
$ws.Range("H11").Value = 25.131568
$ws.Range("I11").Value = 0.002732438555542716
$ws.Range("J11").Value = 0.002732438555542716
$ws.Range("O11").Value = 0.3510414535684271
$ws.Range("P11").Value = 0.3510414535684271
$ws.Range("Q11").Value = 138.1190878467058
$ws.Range("R11").Value = 1243.071790620352
$ws.Range("S11").Value = 0.0009591992023241283
$ws.Range("T11").Value = 0.0009591992023241284
# Row 12
$ws.Range("G12").Value = 8.377189333333332
$ws.Range("H12").Value = 25.131568
$ws.Range("I12").Value = 0.002732438555542716
$ws.Range("J12").Value = 0.002732438555542716
$ws.Range("M12").Value = 5.464566666666666
$ws.Range("N12").Value = 16.3937
$ws.Range("O12").Value = 0.1163479571613943
$ws.Range("P12").Value = 0.1163479571613943
$ws.Range("Q12").Value = 45.77770959128888
$ws.Range("R12").Value = 411.9993863216
$ws.Range("S12").Value = 0.0003179136440064259
$ws.Range("T12").Value = 0.000317913644006426
# Row 13
$ws.Range("G13").Value = 8.377189333333332
$ws.Range("H13").Value = 25.131568
$ws.Range("I13").Value = 0.002732438555542716
$ws.Range("J13").Value = 0.002732438555542716
$ws.Range("M13").Value = 17.60745533333333
$ws.Range("N13").Value = 52.822366
$ws.Range("O13").Value = 0.3748863512527063
$ws.Range("P13").Value = 0.3748863512527063
$ws.Range("Q13").Value = 147.5009870055431
$ws.Range("R13").Value = 1327.508883049888
$ws.Range("S13").Value = 0.001024353920109624
$ws.Range("T13").Value = 0.001024353920109624
# Row 14
$ws.Range("G14").Value = 2.489778666666667
$ws.Range("H14").Value = 7.469336
$ws.Range("I14").Value = 0.0008121061793956991
$ws.Range("J14").Value = 0.0008121061793956991
$ws.Range("M14").Value = 7.407905
$ws.Range("N14").Value = 22.223715
$ws.Range("O14").Value = 0.1577242380174723
$ws.Range("P14").Value = 0.1577242380174723
$ws.Range("Q14").Value = 18.44404383369333
$ws.Range("R14").Value = 165.99639450324
$ws.Range("S14").Value = 0.0001280888283344673
$ws.Range("T14").Value = 0.0001280888283344673
# Row 15
$ws.Range("G15").Value = 2.489778666666667
$ws.Range("H15").Value = 7.469336
$ws.Range("I15").Value = 0.0008121061793956991
$ws.Range("J15").Value = 0.0008121061793956991
$ws.Range("O15").Value = 0.3510414535684271
$ws.Range("P15").Value = 0.3510414535684271
$ws.Range("Q15").Value = 41.05027888194489
$ws.Range("R15").Value = 369.452509937504
$ws.Range("S15").Value = 0.000285082933666968
$ws.Range("T15").Value = 0.000285082933666968
# Row 16
$ws.Range("G16").Value = 2.489778666666667
$ws.Range("H16").Value = 7.469336
$ws.Range("I16").Value = 0.0008121061793956991
$ws.Range("J16").Value = 0.0008121061793956991
$ws.Range("M16").Value = 5.464566666666666
$ws.Range("N16").Value = 16.3937
$ws.Range("O16").Value = 0.1163479571613943
$ws.Range("P16").Value = 0.1163479571613943
$ws.Range("Q16").Value = 13.60556150924444
$ws.Range("R16").Value = 122.4500535832
$ws.Range("S16").Value = 0.00009448689497083437
$ws.Range("T16").Value = 0.00009448689497083438
# Row 17
$ws.Range("G17").Value = 2.489778666666667
$ws.Range("H17").Value = 7.469336
$ws.Range("I17").Value = 0.0008121061793956991
$ws.Range("J17").Value = 0.0008121061793956991
$ws.Range("M17").Value = 17.60745533333333
$ws.Range("N17").Value = 52.822366
$ws.Range("O17").Value = 0.3748863512527063
$ws.Range("P17").Value = 0.3748863512527063
$ws.Range("Q17").Value = 43.83866666321956
$ws.Range("R17").Value = 394.547999968976
$ws.Range("S17").Value = 0.0003044475224234293
$ws.Range("T17").Value = 0.0003044475224234294
